$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Set text labels first so shared-string table indices match the
#     order the strings were typed in originally: Parciales, Final,
#     Tarea, Definitivas, Computacionales ---
$ws.Range("H13").Value = "Parciales"
$ws.Range("I13").Value = "Final"
$ws.Range("J13").Value = "Tarea"
$ws.Range("B8").Value = "Definitivas"
$ws.Range("O12").Value = "Computacionales"

# Row 8
$ws.Range("C8").Value = 4.3

# Row 9
$ws.Range("C9").Value = 4.7

# Row 10
$ws.Range("C10").Value = 4.2

# Row 11
$ws.Range("C11").Value = 3.3

# Row 12
$ws.Range("C12").Value = 3.4

# Row 13
$ws.Range("C13").Formula = "=AVERAGE(E18:G18)"

# Row 14
$ws.Range("A14").Value = 4.07
$ws.Range("B14").Value = 75
$ws.Range("C14").Formula = "=AVERAGE(C8:C13)"
$ws.Range("D14").Value = 18
$ws.Range("H14").Value = 0.6
$ws.Range("I14").Value = 0.25
$ws.Range("J14").Value = 0.15
$ws.Range("K14").Formula = "=SUM(H14:J14)"
$ws.Range("O14").Value = 85
$ws.Range("P14").Formula = "=18/25"

# Row 15
$ws.Range("A15").Formula = "=A14*B14"
$ws.Range("B15").Formula = "=D14*C14"
$ws.Range("C15").Formula = "=(B15+A15)/(D14+B14)"
$ws.Range("H15").Value = 4.2
$ws.Range("I15").Formula = "=24/29*5"
$ws.Range("J15").Value = 5
$ws.Range("O15").Value = 100
$ws.Range("P15").Formula = "=P14*5"

# Row 16
$ws.Range("H16").Formula = "=H15*H14"
$ws.Range("I16").Formula = "=I15*I14"
$ws.Range("J16").Formula = "=J15*J14"
$ws.Range("K16").Formula = "=SUM(H16:J16)"
$ws.Range("O16").Value = 100

# Row 17
$ws.Range("O17").Value = 93

# Row 18
$ws.Range("E18").Value = 4.5999999999999996
$ws.Range("F18").Value = 3.4
$ws.Range("G18").Value = 3.6
$ws.Range("O18").Value = 90

# Row 19
$ws.Range("N19").Value = 1
$ws.Range("O19").Formula = "=AVERAGE(O13:O18)*5/100"

# Row 20
$ws.Range("O20").Formula = "=O19+0.1*N19"

# Column B width (character width ~11.7109375; closest value reachable
# through the ColumnWidth COM setter's internal pixel quantization)
$ws.Range("B1").ColumnWidth = 10.8

# Sheet view adjustments - move the selection to C9 (matches the final
# cursor position after entering the new "Computacionales" assignment
# grades); the workbook/window scroll-position chrome (topLeftCell,
# xWindow/yWindow, absPath) is local UI state that isn't reachable
# through the exposed COM surface here, so it is left as-is.
$ws.Range("C9").Select()
